$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 1500
$ws.Range("I10").Value = 1500
$ws.Range("K10").Value = 1500
$ws.Range("M10").Value = -1207
$ws.Range("H28").Value = 1541.5
$ws.Range("I28").Value = 1755.5555
$ws.Range("J28").Value = 899.3333
$ws.Range("K28").Value = 1755.5555
$ws.Range("L28").Value = 899.3333
$ws.Range("M28").Value = -1270.5555
$ws.Range("N28").Value = -1869.3333
$ws.Range("H41").Value = 1218.1666
$ws.Range("I41").Value = 211.25
$ws.Range("K41").Value = 211.25
$ws.Range("M41").Value = 228.75
$ws.Range("H62").Value = 7291.933
$ws.Range("I62").Value = 6391.9
$ws.Range("K62").Value = 6391.9
$ws.Range("M62").Value = -5767.9
$ws.Range("H65").Value = 7291.933
$ws.Range("I65").Value = 6391.9
$ws.Range("K65").Value = 31959.5
$ws.Range("M65").Value = -28839.5
$ws.Range("H98").Value = 3174.8462
$ws.Range("I98").Value = 2002.35
$ws.Range("J98").Value = 7083.1665
$ws.Range("K98").Value = 2002.35
$ws.Range("L98").Value = 7083.1665
$ws.Range("M98").Value = -504.3499999999999
$ws.Range("N98").Value = -10079.1665
$ws.Range("H122").Value = 3174.8462
$ws.Range("I122").Value = 2002.35
$ws.Range("J122").Value = 7083.1665
$ws.Range("K122").Value = 6007.049999999999
$ws.Range("L122").Value = 21249.4995
$ws.Range("M122").Value = -3557.049999999999
$ws.Range("N122").Value = -26149.4995
$ws.Range("H132").Value = 8405508
$ws.Range("I132").Value = 8405508
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 25216524
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -25213994
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 1419856.6
$ws.Range("I138").Value = 2290.889
$ws.Range("J138").Value = 2109483
$ws.Range("K138").Value = 6872.667
$ws.Range("L138").Value = 6328449
$ws.Range("M138").Value = -1732.667
$ws.Range("N138").Value = -6338729

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18341.967
$ws.Range("I32").Value = 20452.129
$ws.Range("K32").Value = 20452.129
$ws.Range("M32").Value = -20165.129
$ws.Range("H45").Value = 5330
$ws.Range("I45").Value = 4643.6
$ws.Range("K45").Value = 4643.6
$ws.Range("M45").Value = -4266.6
$ws.Range("H74").Value = 3734.28
$ws.Range("I74").Value = 1500.5312
$ws.Range("K74").Value = 1500.5312
$ws.Range("M74").Value = -626.5311999999999
$ws.Range("H77").Value = 3734.28
$ws.Range("I77").Value = 1500.5312
$ws.Range("K77").Value = 7502.655999999999
$ws.Range("M77").Value = -3134.655999999999
$ws.Range("H132").Value = 1604.2642
$ws.Range("I132").Value = 1352.9767
$ws.Range("K132").Value = 4058.9301
$ws.Range("M132").Value = -1528.9301

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 297.58334
$ws.Range("I22").Value = 297.58334
$ws.Range("K22").Value = 297.58334
$ws.Range("M22").Value = -124.58334
$ws.Range("H94").Value = 5883211
$ws.Range("I94").Value = 999.3333
$ws.Range("K94").Value = 999.3333
$ws.Range("M94").Value = -548.3333
$ws.Range("H134").Value = 12703.679
$ws.Range("I134").Value = 13334.962
$ws.Range("K134").Value = 40004.886
$ws.Range("M134").Value = -37469.886

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8797.4
$ws.Range("J99").Value = 8197.4
$ws.Range("L99").Value = 8197.4
$ws.Range("N99").Value = -11193.4
$ws.Range("H126").Value = 8797.4
$ws.Range("J126").Value = 8197.4
$ws.Range("L126").Value = 24592.2
$ws.Range("N126").Value = -29532.2
$ws.Range("H132").Value = 2139.1333
$ws.Range("I132").Value = 1692.0625
$ws.Range("J132").Value = 2650.0715
$ws.Range("K132").Value = 5076.1875
$ws.Range("L132").Value = 7950.2145
$ws.Range("M132").Value = -2546.1875
$ws.Range("N132").Value = -13010.2145
$ws.Range("H134").Value = 2560.3928
$ws.Range("I134").Value = 1858.409
$ws.Range("K134").Value = 5575.227000000001
$ws.Range("M134").Value = -3040.227000000001
$ws.Range("H135").Value = 99998.5
$ws.Range("J135").Value = 99998.5
$ws.Range("L135").Value = 99998.5
$ws.Range("N135").Value = -110138.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1560.2911
$ws.Range("J107").Value = 2047.3673
$ws.Range("L107").Value = 6142.1019
$ws.Range("N107").Value = -9982.1019
$ws.Range("H121").Value = 565691.2
$ws.Range("J121").Value = 793088.5600000001
$ws.Range("L121").Value = 2379265.68
$ws.Range("N121").Value = -2381885.68
$ws.Range("H131").Value = 2370.9167
$ws.Range("I131").Value = 3707.25
$ws.Range("J131").Value = 1925.4722
$ws.Range("K131").Value = 11121.75
$ws.Range("L131").Value = 5776.4166
$ws.Range("M131").Value = -6081.75
$ws.Range("N131").Value = -15856.4166

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 711.88464
$ws.Range("I97").Value = 539.7857
$ws.Range("J97").Value = 912.6667
$ws.Range("K97").Value = 539.7857
$ws.Range("L97").Value = 912.6667
$ws.Range("M97").Value = -43.78570000000002
$ws.Range("N97").Value = -1904.6667
$ws.Range("H132").Value = 2741.6086
$ws.Range("I132").Value = 2566.5293
$ws.Range("J132").Value = 3237.6667
$ws.Range("K132").Value = 7699.5879
$ws.Range("L132").Value = 9713.000100000001
$ws.Range("M132").Value = -5169.5879
$ws.Range("N132").Value = -14773.0001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1774.75
$ws.Range("H40").Value = 4134.1113
$ws.Range("I40").Value = 3910.5
$ws.Range("J40").Value = 4581.3335
$ws.Range("K40").Value = 3910.5
$ws.Range("L40").Value = 4581.3335
$ws.Range("M40").Value = -3774.5
$ws.Range("N40").Value = -4853.3335
$ws.Range("H55").Value = 1842.95
$ws.Range("J55").Value = 2251.5454
$ws.Range("L55").Value = 2251.5454
$ws.Range("N55").Value = -2597.5454
$ws.Range("H132").Value = 5954.2964
$ws.Range("I132").Value = 6581.091
$ws.Range("K132").Value = 19743.273
$ws.Range("M132").Value = -17213.273
$ws.Range("H141").Value = 182857.5
$ws.Range("J141").Value = 182857.5
$ws.Range("L141").Value = 182857.5
$ws.Range("N141").Value = -193217.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5869.115
$ws.Range("I81").Value = 5723.6665
$ws.Range("J81").Value = 6480
$ws.Range("K81").Value = 11447.333
$ws.Range("L81").Value = 12960
$ws.Range("M81").Value = -10386.333
$ws.Range("N81").Value = -15082
$ws.Range("H84").Value = 5869.115
$ws.Range("I84").Value = 5723.6665
$ws.Range("J84").Value = 6480
$ws.Range("K84").Value = 57236.665
$ws.Range("L84").Value = 64800
$ws.Range("M84").Value = -51932.665
$ws.Range("N84").Value = -75408
$ws.Range("H122").Value = 2862.6924
$ws.Range("I122").Value = 3014.6382
$ws.Range("K122").Value = 9043.9146
$ws.Range("M122").Value = -6593.9146
$ws.Range("H132").Value = 18930.48
$ws.Range("I132").Value = 23324.846
$ws.Range("K132").Value = 69974.538
$ws.Range("M132").Value = -67444.538
$ws.Range("H136").Value = 8912.23
$ws.Range("J136").Value = 5117.4707
$ws.Range("L136").Value = 15352.4121
$ws.Range("N136").Value = -20452.4121
